# Apply the "lines_states" update: two new line rows (line7, line8) are
# inserted into the table (after line6), pushing the existing extr1..extr8
# rows down by two rows, and a number of from_bus/to_bus/in_service values
# change as part of the re-run (contingencies with rene fine).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the row-1/A-column style (bold, centered, bordered "header/index"
# style) down onto the two brand-new rows (16 and 17) by copying the format
# from the row directly above before we overwrite the values.
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(17, 1))

# Final target values for rows 8..17, columns A (index), B (name),
# C (from_bus), D (to_bus), E (in_service).
$data = @(
    @(8,  6,  "line7", 14, 11, $false),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $true),
    @(16, 14, "extr7", 5,  7,  $false),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
    $r   = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
